# Update "want-to-go" head counts (column F) on several sheets.
# Sheet layout: 1=展览 (Exhibition), 2=演出 (Performance), 3=本地生活 (Local life), 4=全部类型 (All types)

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheetId 1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value  = 736
$ws1.Range("F13").Value = 255
$ws1.Range("F16").Value = 1291
$ws1.Range("F19").Value = 1086
$ws1.Range("F21").Value = 1282
$ws1.Range("F22").Value = 649
$ws1.Range("F24").Value = 1242
$ws1.Range("F28").Value = 930
$ws1.Range("F29").Value = 15
$ws1.Range("F31").Value = 1319

# --- Sheet "演出" (sheetId 2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 6

# --- Sheet "全部类型" (sheetId 4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F14").Value = 6
$ws4.Range("F16").Value = 736
$ws4.Range("F26").Value = 255
$ws4.Range("F29").Value = 1291
$ws4.Range("F32").Value = 1086
$ws4.Range("F34").Value = 1282
$ws4.Range("F35").Value = 649
$ws4.Range("F37").Value = 1242
$ws4.Range("F43").Value = 930
$ws4.Range("F44").Value = 15
$ws4.Range("F46").Value = 1319

$wb.Save()
